$d = $word.ActiveDocument

$find = "Сазвежђе сазвежђе Орион 2022: 16-25 јануар, 14-23 фебруар, 14-24 март"
$replace = "Сазвежђе сазвежђе Орион током 2022. године посматрамо 16-25 јануар, 14-23 фебруар, 14-24 март"

$range = $d.Content
$range.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
